# Update the "Förändrad" (changed) date column from 2023-10-13 (45212)
# to 2023-10-22 (45221) for rows 2 through 16 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
